# Updated cryptos list (price/volume refresh) on Thu Apr 11 11:27:29 UTC 2024 with GitHub Actions.
# Cells in column D hold plain-text, locale-formatted numbers (e.g. "598.23", "70.569.76").
# A leading apostrophe forces Excel to store the assigned value as text (quote-prefixed,
# General number format) instead of silently parsing it into a numeric cell value, which
# keeps the written cells equivalent to the original inlineStr text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.569.76'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').Value = '3.561.38'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('D5').Value = '''598.23'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').Value = '''172.44'
$ws.Range('E6').Value = '  +1.48%  '
$ws.Range('D7').Value = '3.553.76'
$ws.Range('E7').Value = '  +1.26%  '
$ws.Range('D8').Value = '''0.613'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  +4.59%  '
$ws.Range('D11').Value = '''7.43'
$ws.Range('E11').Value = '  +9.28%  '
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '''46.40'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').Value = '4.133.43'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').Value = '''611.68'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').Value = '3.562.46'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '70.603.76'
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('D21').Value = '''17.35'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').Value = '''9.24'
$ws.Range('E23').Value = '  -16.63%  '
$ws.Range('D24').Value = '''15.70'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').Value = '''96.73'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('D26').Value = '''3.72'
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '''2.61'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '''33.85'
$ws.Range('E29').Value = '  +3.95%  '
$ws.Range('D30').Value = '''9.05'
$ws.Range('E30').Value = '  -1.57%  '
$ws.Range('D31').Value = '''8.27'
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('D33').Value = '''667.04'
$ws.Range('E33').Value = '  +9.27%  '
$ws.Range('D34').Value = '''7.13'
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('D35').Value = '''1.30'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').Value = '''3.62'
$ws.Range('E36').Value = '  +4.43%  '
$ws.Range('D37').Value = '''0.101'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0475'
$ws.Range('E39').Value = '  +7.28%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''57.40'
$ws.Range('E40').Value = '  +0.54%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = '''0.142'
$ws.Range('E42').Value = '  +4.82%  '
$ws.Range('D43').Value = '3.378.29'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('E44').Value = '  -1.93%  '
$ws.Range('D45').Value = '0.0₃0704'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('D46').Value = '''32.73'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').Value = '''2.94'
$ws.Range('E47').Value = '  +7.19%  '
$ws.Range('E48').Value = '  +4.04%  '
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').Value = '''132.35'
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('E51').Value = '  -0.08%  '
